# The deck ships two theme parts:
#   ppt/theme/theme1.xml  ("Integral"     colours) -> used by the slide master
#   ppt/theme/theme2.xml  ("Office Theme" colours) -> used by the notes master
#
# The target revision swaps the theme contents between the two parts (the
# slide master ends up with the "Office Theme" palette, the notes master
# ends up with the "Integral" palette) while everything else - fonts,
# format scheme, part relationships - stays untouched (those two blocks
# were already byte-identical between the parts).
#
# PowerPoint's automation surface lets us rewrite a theme's 12-slot colour
# scheme via Slide.ThemeColorScheme (dk1, lt1, dk2, lt2, accent1-6, hlink,
# folHlink, in that order) - this writes straight into the clrScheme of the
# theme part backing the slide master (theme1.xml here).

$p = $ppt.ActivePresentation

function Set-ThemeColors($themeColorScheme, [string[]]$hexColors) {
    for ($i = 0; $i -lt $hexColors.Count; $i++) {
        $hex = $hexColors[$i]
        $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
        $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
        $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
        # OLE RGB longs are packed 0x00BBGGRR
        $ole = $r + ($g * 256) + ($b * 65536)
        $themeColorScheme.Colors($i + 1).RGB = $ole
    }
}

# Order: dk1, lt1, dk2, lt2, accent1, accent2, accent3, accent4, accent5,
#        accent6, hlink, folHlink
$officeThemeColors = @(
    "000000", "FFFFFF", "44546A", "E7E6E6",
    "5B9BD5", "ED7D31", "A5A5A5", "FFC000",
    "4472C4", "70AD47", "0563C1", "954F72"
)

$slide = $p.Slides.Item(1)
Set-ThemeColors $slide.ThemeColorScheme $officeThemeColors
